# SponsorSynq Competitor Comparison - "Action Items" sheet maintenance.
#
# Remove the three already-shipped "HIGH PRIORITY" rows (Sponsor Matching,
# Verification System, Dashboard Savings Display - all marked COMPLETED)
# plus the "Ambassador Visibility" row, and let every row below shift up to
# fill the gap. No other content on the sheet changes: the five remaining
# "HIGH" rows (Event Collaboration/Co-hosting, Promoter Referral System,
# Landing Page Overhaul, Venue Partnership System, Revenue Stream
# Documentation) and the whole MEDIUM/LOW priority tables just move up by
# four rows, carrying their existing formatting, row heights and merged
# cells with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Action Items")

# Rows 10-13 (the 3 COMPLETED items + Ambassador Visibility) go away;
# Delete() shifts rows 14-30 up into their place automatically, carrying
# along styles, row heights and merged-cell ranges.
$ws.Rows("10:13").Delete()
